$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Capture existing header texts before we move things around (Value2 avoids
# the broken reflection string that .Value returns on this host).
$distReal = $ws.Range("A1").Value2
$distCalc = $ws.Range("B1").Value2

# Copy the existing header style (s="1") onto the cells that will become
# headers so it's reused rather than a brand-new style entry created.
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the old "Erro" header; re-home "Distância real" / "Distância
# calculada" two columns to the right (A/B -> E/F).
$ws.Range("C1").ClearContents()
$ws.Range("E1").Value = $distReal
$ws.Range("F1").Value = $distCalc

# New header labels for the inserted columns A-D
$ws.Range("A1").Value = "Idade"
$ws.Range("B1").Value = "Altura"
$ws.Range("C1").Value = "Peso"
$ws.Range("D1").Value = "Gênero"

# --- Data row 2 ---------------------------------------------------------
# A2, B2, C2, E2 need to hold literal *text* (e.g. "19"), not numbers, even
# though the text looks numeric. A direct .Value assignment lets the host
# auto-convert numeric-looking strings into real numbers, so we stage the
# text in a scratch cell formatted as Text, then bring only the *value*
# (xlPasteValues) over - this keeps the text typing without leaving any
# stray number-format/style behind on the target cell.
function Set-TextValue($range, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
    $scratch.Clear()
}

Set-TextValue $ws.Range("A2") "19"
Set-TextValue $ws.Range("B2") "1.74"
Set-TextValue $ws.Range("C2") "84"
$ws.Range("D2").Value = "Masculino"
Set-TextValue $ws.Range("E2") "0"
$ws.Range("F2").Value = -0.7542467618611298
